$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "39.902.96"
$ws.Range("E2").Value = "  -0.32%  "

# Row 3
$ws.Range("D3").Value = "2.231.86"
$ws.Range("E3").Value = "  -4.82%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.86"
$ws.Range("E5").Value = "  -5.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.23"
$ws.Range("E6").Value = "  +0.76%  "

# Row 7
$ws.Range("E7").Value = "  -2.25%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -2.29%  "

# Row 10
$ws.Range("E10").Value = "  -0.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.42"
$ws.Range("E11").Value = "  +1.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.83"
$ws.Range("E12").Value = "  -8.87%  "

# Row 13
$ws.Range("E13").Value = "  -2.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.36"
$ws.Range("E14").Value = "  -0.80%  "

# Row 15
$ws.Range("D15").Value = "2.578.20"
$ws.Range("E15").Value = "  -4.72%  "

# Row 16
$ws.Range("E16").Value = "  -4.24%  "

# Row 17
$ws.Range("D17").Value = "2.234.42"
$ws.Range("E17").Value = "  -5.56%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.724"
$ws.Range("E18").Value = "  -4.40%  "

# Row 19
$ws.Range("D19").Value = "39.822.42"
$ws.Range("E19").Value = "  -0.55%  "

# Row 20
$ws.Range("E20").Value = "  -0.96%  "

# Row 21
$ws.Range("E21").Value = "  -4.75%  "

# Row 22
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.61"
$ws.Range("E22").Value = "  +0.41%  "

# Row 23
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.41"
$ws.Range("E23").Value = "  -3.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "232.00"
$ws.Range("E24").Value = "  -1.24%  "

# Row 25
$ws.Range("E25").Value = "  +0.06%  "

# Row 26
$ws.Range("E26").Value = "  -4.85%  "

# Row 27
$ws.Range("E27").Value = "  +1.22%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").Value = "  +4.40%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.96"
$ws.Range("E29").Value = "  -2.54%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.23"
$ws.Range("E30").Value = "  +0.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.21"
$ws.Range("E31").Value = "  +1.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.03"
$ws.Range("E32").Value = "  -4.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.31%  "

# Row 34
$ws.Range("E34").Value = "  -4.75%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0708"
$ws.Range("E35").Value = "  -0.97%  "

# Row 36
$ws.Range("E36").Value = "  -5.14%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.55"
$ws.Range("E37").Value = "  +6.75%  "

# Row 38
$ws.Range("E38").Value = "  -1.58%  "

# Row 39
$ws.Range("E39").Value = "  -0.73%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.68"
$ws.Range("E40").Value = "  -3.38%  "

# Row 41
$ws.Range("E41").Value = "  -2.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.74"
$ws.Range("E42").Value = "  -2.87%  "

# Row 43
$ws.Range("D43").Value = "1.953.81"
$ws.Range("E43").Value = "  -0.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.19"
$ws.Range("E44").Value = "  -2.99%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0270"
$ws.Range("E45").Value = "  +1.91%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.46"
$ws.Range("E46").Value = "  +0.62%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.36"
$ws.Range("E47").Value = "  -6.66%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("E48").Value = "  -2.52%  "

# Row 49
$ws.Range("D49").Value = "2.447.52"
$ws.Range("E49").Value = "  -4.63%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.79"
$ws.Range("E50").Value = "  +0.91%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.45"
$ws.Range("E51").Value = "  +6.67%  "

